# Apply the edits described by the diff:
#  - update the date paragraph
#  - update the division problems inside the table (5 populated rows x 5 cols)

$d = $word.ActiveDocument

# 1. Update the date line.
$d.Content.Find.Execute("2023-11-25 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-11-26 Sunday", 2)

# 2. Update the division problems. The table has 20 rows (every 4th row, i.e.
#    rows 1, 5, 9, 13, 17 in 1-based indexing, holds the 5 problems; the rows
#    in between are blank spacer rows). Because some problem strings repeat
#    across the sheet (e.g. "71÷6=" and "94÷6=" each appear twice) we address
#    cells directly by row/column instead of relying on a global Find/Replace.

$table = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "18÷7=" },
    @{ Row = 1;  Col = 2; New = "11÷7=" },
    @{ Row = 1;  Col = 3; New = "26÷6=" },
    @{ Row = 1;  Col = 4; New = "12÷5=" },
    @{ Row = 1;  Col = 5; New = "37÷3=" },

    @{ Row = 5;  Col = 1; New = "89÷6=" },
    @{ Row = 5;  Col = 2; New = "98÷4=" },
    @{ Row = 5;  Col = 3; New = "93÷7=" },
    @{ Row = 5;  Col = 4; New = "31÷8=" },
    @{ Row = 5;  Col = 5; New = "64÷8=" },

    @{ Row = 9;  Col = 1; New = "60÷3=" },
    @{ Row = 9;  Col = 2; New = "78÷4=" },
    @{ Row = 9;  Col = 3; New = "54÷7=" },
    @{ Row = 9;  Col = 4; New = "87÷5=" },
    @{ Row = 9;  Col = 5; New = "54÷9=" },

    @{ Row = 13; Col = 1; New = "15÷7=" },
    @{ Row = 13; Col = 2; New = "70÷5=" },
    @{ Row = 13; Col = 3; New = "32÷8=" },
    @{ Row = 13; Col = 4; New = "15÷5=" },
    @{ Row = 13; Col = 5; New = "82÷9=" },

    @{ Row = 17; Col = 1; New = "89÷8=" },
    @{ Row = 17; Col = 2; New = "55÷3=" },
    @{ Row = 17; Col = 3; New = "79÷8=" },
    @{ Row = 17; Col = 4; New = "44÷8=" },
    @{ Row = 17; Col = 5; New = "48÷8=" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}

$d.Save()
